$d = $word.ActiveDocument

# --- 1. Update the "compatible with" paragraph ------------------------------------
# Find the two runs that currently read:
#   "This sample is compatible with the " + "Windows 10 Fall Creators Update SDK (16299)"
# and turn them into a single run reading:
#   "This sample is compatible with the Windows 10 April 2018 Update SDK (17134)"
# with updated paragraph/run formatting (pStyle Heading1, spacing before=0,
# rFonts eastAsiaTheme=minorHAnsi/cs=Times New Roman, color auto, sz 20, szCs 22).

$found = $d.Content.Find.Execute(
    "This sample is compatible with the Windows 10 Fall Creators Update SDK (16299)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This sample is compatible with the Windows 10 April 2018 Update SDK (17134)",
    2)

# Locate the paragraph that now holds the replaced text and restyle it.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*This sample is compatible with the Windows 10 April 2018 Update SDK (17134)*") {
        $p.Style = "Heading1"
        $p.Format.SpaceBefore = 0

        $r = $p.Range
        $r.Font.Italic = $true
        $r.Font.Color = -16777216
        $r.Font.Size = 10
        $r.Font.NameFarEast = "+mn-ea"
        break
    }
}

# --- 2. Footer copyright year 2017 -> 2018 ----------------------------------------
foreach ($sec in $d.Sections) {
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            $ftr.Range.Find.Execute("2017", $true, $false, $false, $false, $false, $true, 1, $false, "2018", 2)
        }
    }
}
